$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Normalize rows 39-47 (the 4 existing "recent task" rows plus the
#        5 about to be appended) to the regular row formatting used by every
#        other task row, by copying the format from row 38. This both moves
#        the highlight off of rows 39-42 and formats the freshly appended
#        rows 43-47 consistently. Done one row at a time for a well defined
#        format-paste. ---
for ($row = 39; $row -le 47; $row++) {
    $ws.Range("C38:E38").Copy() | Out-Null
    $ws.Range("C$row`:E$row").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}
$excel.CutCopyMode = 0

# --- 2. Fill in the values/text for the 5 newly added tasks ---
$ws.Range("C43").Value = 20.3
$ws.Range("D43").Value = "Create i Update za employee"
$ws.Range("E43").Value = "DONE"

$ws.Range("C44").Value = 21
$ws.Range("D44").Value = "Loading Spinner"
$ws.Range("E44").Value = "DONE"

$ws.Range("C45").Value = 22
$ws.Range("D45").Value = "BaseDialog for errors"
$ws.Range("E45").Value = "DONE"

$ws.Range("C46").Value = 23
$ws.Range("D46").Value = "Route transition"
$ws.Range("E46").Value = "DONE"

$ws.Range("C47").Value = 24
$ws.Range("D47").Value = "NotFound page"
$ws.Range("E47").Value = "DONE"

# --- 3. Move the active selection down to the next empty cell below the
#        newly appended rows (mirrors the user continuing to work downward). ---
$ws.Activate() | Out-Null
$ws.Range("C48").Select() | Out-Null
